$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers must be forced to Text format
# so Excel keeps them as strings (matching the source data which uses them as labels).
$ws.Range('D2').Value = '45.917.54'
$ws.Range('E2').Value = '  -2.08%  '
$ws.Range('D3').Value = '2.381.45'
$ws.Range('E3').Value = '  +3.19%  '
$ws.Range('E4').Value = '  -0.07%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '300.12'
$ws.Range('E5').Value = '  -1.15%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '98.23'
$ws.Range('E6').Value = '  -3.31%  '
$ws.Range('E7').Value = '  -1.08%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('E9').Value = '  -4.30%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '34.42'
$ws.Range('E10').Value = '  -6.98%  '
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D11').Value = '0.0787'
$ws.Range('E11').Value = '  -2.02%  '
$ws.Range('E12').Value = '  -4.65%  '
$ws.Range('E13').Value = '  -0.37%  '
$ws.Range('D14').Value = '2.748.48'
$ws.Range('E14').Value = '  +3.33%  '
$ws.Range('D15').Value = '2.376.43'
$ws.Range('E15').Value = '  +3.01%  '
$ws.Range('E16').Value = '  +0.50%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '13.75'
$ws.Range('E17').Value = '  -2.00%  '
$ws.Range('D18').Value = '45.872.87'
$ws.Range('E18').Value = '  -2.19%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '12.72'
$ws.Range('E19').Value = '  -6.83%  '
$ws.Range('E20').Value = '  -0.27%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '6.04'
$ws.Range('E21').Value = '  -1.39%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '66.76'
$ws.Range('E22').Value = '  -0.23%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '243.35'
$ws.Range('E23').Value = '  -2.58%  '
$ws.Range('E24').Value = '  -5.71%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '0.999'
$ws.Range('E25').Value = '  +0.09%  '
$ws.Range('E26').Value = '  -2.51%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '39.06'
$ws.Range('E27').Value = '  -11.99%  '
$ws.Range('E28').Value = '  -3.31%  '
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '9.71'
$ws.Range('E29').Value = '  -2.60%  '
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '21.03'
$ws.Range('E30').Value = '  +3.88%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.75'
$ws.Range('E31').Value = '  +17.34%  '
$ws.Range('E32').Value = '  +4.16%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '5.54'
$ws.Range('E33').Value = '  -4.82%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '147.41'
$ws.Range('E34').Value = '  -0.59%  '
$ws.Range('E35').Value = '  -5.13%  '
$ws.Range('E36').Value = '  -1.32%  '
$ws.Range('E37').Value = '  +5.69%  '
$ws.Range('E38').Value = '  -3.13%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '14.87'
$ws.Range('E39').Value = '  -8.24%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '3.87'
$ws.Range('E40').Value = '  -4.40%  '
$ws.Range('E41').Value = '  -2.98%  '
$ws.Range('E42').Value = '  -7.91%  '
$ws.Range('D43').Value = '1.946.21'
$ws.Range('E43').Value = '  +4.51%  '
$ws.Range('E44').Value = '  -0.01%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '94.94'
$ws.Range('E45').Value = '  +6.65%  '
$ws.Range('E46').Value = '  -10.41%  '
$ws.Range('E47').Value = '  +5.02%  '
$ws.Range('E48').Value = '  -6.04%  '
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '99.11'
$ws.Range('E49').Value = '  +1.78%  '
$ws.Range('D50').Value = '2.618.24'
$ws.Range('E50').Value = '  +3.22%  '
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '68.52'
$ws.Range('E51').Value = '  -8.65%  '
